# Auto-generated: apply scheduled-runner value refresh to Seraph_Profits sheets
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 99
$ws.Range("I6").Value = 99
$ws.Range("K6").Value = 297
$ws.Range("M6").Value = -185
$ws.Range("H9").Value = 88.28570999999999
$ws.Range("I9").Value = 87.5
$ws.Range("J9").Value = 89.333336
$ws.Range("K9").Value = 87.5
$ws.Range("L9").Value = 89.333336
$ws.Range("M9").Value = 81.5
$ws.Range("N9").Value = -427.333336
$ws.Range("H47").Value = 14042
$ws.Range("J47").Value = 7911
$ws.Range("L47").Value = 7911
$ws.Range("N47").Value = -9855
$ws.Range("H55").Value = 511.1
$ws.Range("I55").Value = 46
$ws.Range("J55").Value = 627.375
$ws.Range("K55").Value = 46
$ws.Range("L55").Value = 627.375
$ws.Range("M55").Value = 168
$ws.Range("N55").Value = -1055.375
$ws.Range("H86").Value = 4335.2856
$ws.Range("J86").Value = 5537.5
$ws.Range("L86").Value = 5537.5
$ws.Range("N86").Value = -7783.5
$ws.Range("H88").Value = 1513.6154
$ws.Range("I88").Value = 712.6
$ws.Range("K88").Value = 712.6
$ws.Range("M88").Value = -306.6
$ws.Range("H89").Value = 4335.2856
$ws.Range("J89").Value = 5537.5
$ws.Range("L89").Value = 27687.5
$ws.Range("N89").Value = -38919.5
$ws.Range("H91").Value = 1513.6154
$ws.Range("I91").Value = 712.6
$ws.Range("K91").Value = 712.6
$ws.Range("M91").Value = 691.4
$ws.Range("H111").Value = 3716.8572
$ws.Range("I111").Value = 4004.6
$ws.Range("J111").Value = 2997.5
$ws.Range("K111").Value = 12013.8
$ws.Range("L111").Value = 8992.5
$ws.Range("M111").Value = -8946.799999999999
$ws.Range("N111").Value = -15126.5
$ws.Range("H132").Value = 1349.1154
$ws.Range("I132").Value = 1349.1154
$ws.Range("K132").Value = 4047.3462
$ws.Range("M132").Value = -1517.3462
$ws.Range("H138").Value = 4127.871
$ws.Range("I138").Value = 1980
$ws.Range("J138").Value = 4446.074
$ws.Range("K138").Value = 5940
$ws.Range("L138").Value = 13338.222
$ws.Range("M138").Value = -800
$ws.Range("N138").Value = -23618.222

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()  # cell removed in target (was present before)
$ws.Range("H74").Value = 2211.9443
$ws.Range("I74").Value = 1377.4375
$ws.Range("K74").Value = 1377.4375
$ws.Range("M74").Value = -503.4375
$ws.Range("H77").Value = 2211.9443
$ws.Range("I77").Value = 1377.4375
$ws.Range("K77").Value = 6887.1875
$ws.Range("M77").Value = -2519.1875
$ws.Range("H122").Value = 479651
$ws.Range("I122").Value = 669544.8
$ws.Range("K122").Value = 2008634.4
$ws.Range("M122").Value = -2006184.4
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()  # cell removed in target (was present before)

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 996.6667
$ws.Range("I20").Value = 996.6667
$ws.Range("K20").Value = 996.6667
$ws.Range("M20").Value = -749.6667
$ws.Range("H86").Value = 1463.5
$ws.Range("I86").Value = 1476.1666
$ws.Range("J86").Value = 1450.8334
$ws.Range("K86").Value = 1476.1666
$ws.Range("L86").Value = 1450.8334
$ws.Range("M86").Value = -353.1666
$ws.Range("N86").Value = -3696.8334
$ws.Range("H89").Value = 1463.5
$ws.Range("I89").Value = 1476.1666
$ws.Range("J89").Value = 1450.8334
$ws.Range("K89").Value = 7380.833000000001
$ws.Range("L89").Value = 7254.166999999999
$ws.Range("M89").Value = -1764.833000000001
$ws.Range("N89").Value = -18486.167
$ws.Range("H99").Value = 3969.8096
$ws.Range("I99").Value = 3813
$ws.Range("J99").Value = 4361.8335
$ws.Range("K99").Value = 3813
$ws.Range("L99").Value = 4361.8335
$ws.Range("M99").Value = -2315
$ws.Range("N99").Value = -7357.8335
$ws.Range("H105").Value = 3928.3914
$ws.Range("I105").Value = 3026.1052
$ws.Range("K105").Value = 3026.1052
$ws.Range("M105").Value = -1279.1052
$ws.Range("H107").Value = 678.8261
$ws.Range("I107").Value = 672.0952
$ws.Range("K107").Value = 672.0952
$ws.Range("M107").Value = 1247.9048
$ws.Range("H122").Value = 259984.5
$ws.Range("I122").Value = 299979
$ws.Range("J122").Value = 219990
$ws.Range("K122").Value = 299979
$ws.Range("L122").Value = 219990
$ws.Range("M122").Value = -295079
$ws.Range("N122").Value = -229790

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6573.3076
$ws.Range("I31").Value = 3977.75
$ws.Range("K31").Value = 3977.75
$ws.Range("M31").Value = -3682.75
$ws.Range("H34").Value = 6573.3076
$ws.Range("I34").Value = 3977.75
$ws.Range("K34").Value = 3977.75
$ws.Range("M34").Value = -3775.75
$ws.Range("H99").Value = 16607.137
$ws.Range("I99").Value = 13987.223
$ws.Range("J99").Value = 18420.924
$ws.Range("K99").Value = 13987.223
$ws.Range("L99").Value = 18420.924
$ws.Range("M99").Value = -12489.223
$ws.Range("N99").Value = -21416.924
$ws.Range("H107").Value = 592.84
$ws.Range("I107").Value = 420.05554
$ws.Range("J107").Value = 1037.1428
$ws.Range("K107").Value = 420.05554
$ws.Range("L107").Value = 1037.1428
$ws.Range("M107").Value = 1499.94446
$ws.Range("N107").Value = -4877.1428
$ws.Range("H126").Value = 16607.137
$ws.Range("I126").Value = 13987.223
$ws.Range("J126").Value = 18420.924
$ws.Range("K126").Value = 41961.669
$ws.Range("L126").Value = 55262.772
$ws.Range("M126").Value = -39491.669
$ws.Range("N126").Value = -60202.772
$ws.Range("H133").Value = 67145.5
$ws.Range("I133").Value = 42295
$ws.Range("K133").Value = 42295
$ws.Range("M133").Value = -39765

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 709.6667
$ws.Range("I14").Value = 709.6667
$ws.Range("K14").Value = 2129.0001
$ws.Range("M14").Value = -1956.0001
$ws.Range("H26").Value = 195.83333
$ws.Range("I26").Value = 169.25
$ws.Range("J26").Value = 249
$ws.Range("K26").Value = 507.75
$ws.Range("L26").Value = 747
$ws.Range("M26").Value = -219.75
$ws.Range("N26").Value = -1323
$ws.Range("H137").Value = 6986.8667
$ws.Range("I137").Value = 4985.2856
$ws.Range("J137").Value = 8738.25
$ws.Range("K137").Value = 14955.8568
$ws.Range("L137").Value = 26214.75
$ws.Range("M137").Value = -9855.856800000001
$ws.Range("N137").Value = -36414.75
$ws.Range("H139").Value = 4172.778
$ws.Range("I139").Value = 2264.9167
$ws.Range("J139").Value = 7988.5
$ws.Range("K139").Value = 6794.750100000001
$ws.Range("L139").Value = 23965.5
$ws.Range("M139").Value = -1654.750100000001
$ws.Range("N139").Value = -34245.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5750
$ws.Range("J70").Value = 5000
$ws.Range("L70").Value = 5000
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 5750
$ws.Range("J73").Value = 5000
$ws.Range("L73").Value = 5000
$ws.Range("N73").Value = -6872
$ws.Range("H123").Value = 23052.777
$ws.Range("J123").Value = 23052.777
$ws.Range("L123").Value = 23052.777
$ws.Range("N123").Value = -27952.777
$ws.Range("H136").Value = 38086.77
$ws.Range("J136").Value = 38086.77
$ws.Range("L136").Value = 114260.31
$ws.Range("N136").Value = -119360.31

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6211.75
$ws.Range("J16").Value = 13000
$ws.Range("L16").Value = 13000
$ws.Range("N16").Value = -13340
$ws.Range("H22").Value = 10853.571
$ws.Range("I22").Value = 3988.5
$ws.Range("J22").Value = 13599.6
$ws.Range("K22").Value = 3988.5
$ws.Range("L22").Value = 13599.6
$ws.Range("M22").Value = -3693.5
$ws.Range("N22").Value = -14189.6
$ws.Range("H27").Value = 10853.571
$ws.Range("I27").Value = 3988.5
$ws.Range("J27").Value = 13599.6
$ws.Range("K27").Value = 3988.5
$ws.Range("L27").Value = 13599.6
$ws.Range("M27").Value = -3881.5
$ws.Range("N27").Value = -13813.6

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10000
$ws.Range("J41").Value = 10000
$ws.Range("L41").Value = 10000
$ws.Range("N41").Value = -10780
$ws.Range("H51").Value = 37519
$ws.Range("I51").Value = 26666.334
$ws.Range("J51").Value = 70077
$ws.Range("K51").Value = 26666.334
$ws.Range("L51").Value = 70077
$ws.Range("M51").Value = -26156.334
$ws.Range("N51").Value = -71097
$ws.Range("H60").Value = 40000
$ws.Range("J60").Value = 40000
$ws.Range("L60").Value = 40000
$ws.Range("N60").Value = -41644
$ws.Range("H70").Value = 57797.5
$ws.Range("H73").Value = 57797.5
$ws.Range("H81").Value = 15200
$ws.Range("J81").Value = 17500
$ws.Range("L81").Value = 35000
$ws.Range("N81").Value = -37122
$ws.Range("H84").Value = 15200
$ws.Range("J84").Value = 17500
$ws.Range("L84").Value = 175000
$ws.Range("N84").Value = -185608
$ws.Range("H122").Value = 14294071
$ws.Range("I122").Value = 14294071
$ws.Range("K122").Value = 42882213
$ws.Range("M122").Value = -42879763
$ws.Range("H126").Value = 7501.8
$ws.Range("J126").Value = 7503.3335
$ws.Range("L126").Value = 22510.0005
$ws.Range("N126").Value = -27450.0005
$ws.Range("H132").Value = 2382.625
$ws.Range("I132").Value = 2010.3334
$ws.Range("K132").Value = 6031.0002
$ws.Range("M132").Value = -3501.0002
$ws.Range("H136").Value = 5411.778
$ws.Range("I136").Value = 1987.6666
$ws.Range("K136").Value = 5962.9998
$ws.Range("M136").Value = -3412.9998
